$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Food / 120 / 10-31-2025
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Food"
$ws.Range("C6").Value = 120
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "10-31-2025"
$ws.Range("D6").Style = "Normal"

# Row 7: Water / 1000 / 11-01-2025
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Water"
$ws.Range("C7").Value = 1000
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "11-01-2025"
$ws.Range("D7").Style = "Normal"

# Row 8: Water / 1000 / 11-01-2025
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Water"
$ws.Range("C8").Value = 1000
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "11-01-2025"
$ws.Range("D8").Style = "Normal"
